$d = $word.ActiveDocument

# The document currently contains a single paragraph. Clear it out and
# retype the full replacement content (10 new paragraphs + the original
# paragraph's unchanged text) as one block, using `r for paragraph breaks.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.Delete()

$fullText = "Although the phylogeny of extant salamanders may be inferred using molecular data, the inclusion of fossils is more difficult. Fossil salamanders were thus placed on the extant phylogeny using relationships proposed by the most recent phylogenetic assessment of those fossil species (see SOM). However, two major clades of fossil salamander, the scapherpetonids and batrachosauroidids, are problematic, and the placement of these clades remains uncertain. To handle this, we ran analyses over seven different topologies, each with a different placement for those two clades. One topology, where scapherpetonids and batrachosauroidids are sister to each other and outside of Urodela, produced notably different estimates for rates through time (a greatly enhanced spike in rates at the KPg boundary).`r`r1:10 – Scap sister to all, Batrach sister to Crypts (S,((B,Crypt),(Others)))`r11:20- Scap sister to Crypts, Batrach sister to all (B,((S,Crypt),Others)))`r21:30—Scap outgroup to crown, Batrach outgroup to Scap+Crown (B,(S,Crown))`r31:40—Scap and Batrach sister to each other, outgroup to crown ((B,S),(Crown))`r41:50—Scap sister to Batrach + Crown (S,(B,Crown))`r51-60—Scap sister to crown, Batrach sister to Hyno+Crypt (S,((B,(Hyn,Crypt)),Others)`r61-70—Batrach sister to crown, Scap sister to Hyno+Crypt (B,((S,(Hyn,Crypt)),Others)`r`rWe modelled neoteny using the liability model proposed by Felsenstein (#CITE) and expanded by Revell (#CITE). This model treats the multiple observable discrete states as manifestations of a continuously-varying “hidden” or latent trait that evolves by Brownian Motion. This approach has two useful statistical properties. First, although movement of the continuous trait is constant and symmetrical, transition rates between character states can be asymmetrical without having to optimize a large rate matrix. Second, although the rate of change for the continuous trait is constant across the tree, and the thresholds for discrete states are fixed, the actual frequency of transition between states among subclades can vary. This means that some clades can frequently shift between states while other clades may be “locked” in a single state, again without having to estimate transition points along the phylogeny. Biologically, this model is also very appropriate for neoteny in salamanders, as the discrete states described above are known to be influenced by continuously-varying levels of thyroid hormones (#Rose 1996; #CITE). Paedomorphosis in salamanders may occur in anything from the failure of skull bones like the pterygoid to be resorbed to the retention of gills and juvenile coloration at sexual maturity. For the purposes of this study, we discretized this array of forms with four states: direct development (no larval stage), biphasic/metamorphosing, facultative neoteny (some populations metamorphose in certain conditions), and obligate neoteny (no populations metamorphose). "
$d.Paragraphs(1).Range.Text = $fullText

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)

# Relocate the hidden "_GoBack" bookmark: it used to sit between "retention
# of " and "gills" in the (now final) paragraph; it belongs between "Hyn,"
# and "Crypt" in the "51-60" topology paragraph instead.
$anchor = $d.Content
$found = $anchor.Find.Execute("Hyn,", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($anchor.End, $anchor.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    Write-Output "Bookmark _GoBack relocated."
} else {
    Write-Output "WARNING: anchor text 'Hyn,' not found; bookmark not moved."
}
